$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a text-looking value to a cell while forcing text storage
# (prevents Excel from auto-converting values like "302.86" into a number).
function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = "23.731.68"
$ws.Range("E2").Value = "  +0.94%  "
$ws.Range("D3").Value = "1.658.26"
$ws.Range("E3").Value = "  +0.92%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("E5").Value = "  +0.18%  "
Set-TextValue $ws.Range("D6") "302.86"
$ws.Range("E6").Value = "  -0.30%  "
Set-TextValue $ws.Range("D7") "0.3816"
$ws.Range("E7").Value = "  +0.46%  "
Set-TextValue $ws.Range("D8") "0.3612"
$ws.Range("E8").Value = "  -0.27%  "
Set-TextValue $ws.Range("D9") "51.19"
$ws.Range("E9").Value = "  -1.76%  "
Set-TextValue $ws.Range("D10") "0.08199"
$ws.Range("E10").Value = "  +0.14%  "
Set-TextValue $ws.Range("D11") "1.235"
$ws.Range("E11").Value = "  -0.29%  "
$ws.Range("E12").Value = "  +0.16%  "
Set-TextValue $ws.Range("D13") "22.56"
$ws.Range("E13").Value = "  -0.10%  "
Set-TextValue $ws.Range("D14") "6.455"
$ws.Range("E14").Value = "  -0.37%  "
Set-TextValue $ws.Range("D15") "7.441"
$ws.Range("E15").Value = "  +0.85%  "
$ws.Range("E16").Value = "  -1.09%  "
$ws.Range("D17").Value = "1.654.42"
$ws.Range("E17").Value = "  +0.97%  "
Set-TextValue $ws.Range("D18") "97.71"
$ws.Range("E18").Value = "  +2.59%  "
Set-TextValue $ws.Range("D19") "0.07027"
$ws.Range("E19").Value = "  +1.12%  "
Set-TextValue $ws.Range("D20") "6.832"
$ws.Range("E20").Value = "  +3.82%  "
Set-TextValue $ws.Range("D21") "17.62"
$ws.Range("E21").Value = "  +0.07%  "
Set-TextValue $ws.Range("D22") "1.001"
$ws.Range("E22").Value = "  +0.21%  "
Set-TextValue $ws.Range("D23") "12.77"
$ws.Range("E23").Value = "  +1.82%  "
$ws.Range("D24").Value = "23.742.42"
$ws.Range("E24").Value = "  +0.96%  "
Set-TextValue $ws.Range("D25") "2.516"
$ws.Range("E25").Value = "  -0.22%  "
Set-TextValue $ws.Range("D26") "3.007"
$ws.Range("E26").Value = "  -1.99%  "
Set-TextValue $ws.Range("D27") "21.25"
$ws.Range("E27").Value = "  +0.18%  "
Set-TextValue $ws.Range("D28") "153.60"
$ws.Range("E28").Value = "  +1.11%  "
Set-TextValue $ws.Range("D29") "5.227"
$ws.Range("E29").Value = "  -0.67%  "
Set-TextValue $ws.Range("D30") "134.18"
$ws.Range("E30").Value = "  +0.59%  "
$ws.Range("D31").Value = "1.844.53"
$ws.Range("E31").Value = "  +1.41%  "
Set-TextValue $ws.Range("D32") "7.136"
$ws.Range("E32").Value = "  +7.75%  "
Set-TextValue $ws.Range("D33") "2.239"
$ws.Range("E33").Value = "  +4.09%  "
Set-TextValue $ws.Range("D34") "12.06"
$ws.Range("E34").Value = "  +4.95%  "
Set-TextValue $ws.Range("D35") "1.056"
$ws.Range("E35").Value = "  -4.07%  "
Set-TextValue $ws.Range("D36") "0.02818"
$ws.Range("E36").Value = "  +1.72%  "
Set-TextValue $ws.Range("D37") "0.2521"
$ws.Range("E37").Value = "  +0.21%  "
Set-TextValue $ws.Range("D38") "6.117"
$ws.Range("E38").Value = "  +1.59%  "
Set-TextValue $ws.Range("D39") "0.08813"
$ws.Range("E39").Value = "  +0.56%  "
Set-TextValue $ws.Range("D40") "0.06998"
$ws.Range("E40").Value = "  -0.75%  "
Set-TextValue $ws.Range("D41") "12.99"
$ws.Range("E41").Value = "  +5.51%  "
Set-TextValue $ws.Range("D42") "0.7002"
$ws.Range("E42").Value = "  -1.05%  "
Set-TextValue $ws.Range("D43") "1.338"
$ws.Range("E43").Value = "  -1.15%  "
Set-TextValue $ws.Range("D44") "16.11"
$ws.Range("E44").Value = "  +3.23%  "
Set-TextValue $ws.Range("D45") "0.6527"
$ws.Range("E45").Value = "  -0.57%  "
Set-TextValue $ws.Range("D46") "1.000"
$ws.Range("E46").Value = "  +0.18%  "
$ws.Range("E47").Value = "  +0.19%  "
Set-TextValue $ws.Range("D48") "3.967"
$ws.Range("E48").Value = "  -0.10%  "
Set-TextValue $ws.Range("D49") "0.07915"
$ws.Range("E49").Value = "  -0.93%  "
Set-TextValue $ws.Range("D50") "128.49"
$ws.Range("E50").Value = "  -0.38%  "
$ws.Range("E51").Value = "  -0.99%  "
